# Update team-specific time-allocation matrix values on Sheet1 (commit:
# "added team spec time commit pt2"). Each assignment below corresponds to
# one changed <c> value cell from the source OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.196078431372549
$ws.Cells.Item(2, 3).Value = 0.5343137254901961
$ws.Cells.Item(2, 10).Value = 0.02941176470588235
$ws.Cells.Item(2, 16).Value = 0.1617647058823529
$ws.Cells.Item(2, 19).Value = 0.07843137254901961
$ws.Cells.Item(3, 3).Value = 0.03571428571428571
$ws.Cells.Item(3, 10).Value = 0.008928571428571428
$ws.Cells.Item(3, 16).Value = 0.7946428571428571
$ws.Cells.Item(3, 19).Value = 0.1607142857142857
$ws.Cells.Item(4, 10).Value = 0.04878048780487805
$ws.Cells.Item(4, 16).Value = 0.6097560975609756
$ws.Cells.Item(4, 19).Value = 0.3414634146341464
$ws.Cells.Item(6, 2).Value = 0.05579399141630902
$ws.Cells.Item(6, 4).Value = 0.01716738197424893
$ws.Cells.Item(6, 6).Value = 0.06437768240343347
$ws.Cells.Item(6, 10).Value = 0.2446351931330472
$ws.Cells.Item(6, 15).Value = 0.02145922746781116
$ws.Cells.Item(6, 17).Value = 0.1630901287553648
$ws.Cells.Item(6, 18).Value = 0.07296137339055794
$ws.Cells.Item(6, 19).Value = 0.3605150214592275
$ws.Cells.Item(7, 2).Value = 0.08333333333333333
$ws.Cells.Item(7, 4).Value = 0.009803921568627451
$ws.Cells.Item(7, 6).Value = 0.04411764705882353
$ws.Cells.Item(7, 10).Value = 0.09313725490196079
$ws.Cells.Item(7, 15).Value = 0.03431372549019608
$ws.Cells.Item(7, 17).Value = 0.1911764705882353
$ws.Cells.Item(7, 18).Value = 0.07352941176470588
$ws.Cells.Item(7, 19).Value = 0.4705882352941176
$ws.Cells.Item(8, 2).Value = 0.06790123456790123
$ws.Cells.Item(8, 4).Value = 0.01646090534979424
$ws.Cells.Item(8, 5).Value = 0.00205761316872428
$ws.Cells.Item(8, 6).Value = 0.05349794238683128
$ws.Cells.Item(8, 10).Value = 0.1193415637860082
$ws.Cells.Item(8, 15).Value = 0.01646090534979424
$ws.Cells.Item(8, 17).Value = 0.2160493827160494
$ws.Cells.Item(8, 18).Value = 0.1172839506172839
$ws.Cells.Item(8, 19).Value = 0.3909465020576132
$ws.Cells.Item(9, 2).Value = 0.04705882352941176
$ws.Cells.Item(9, 6).Value = 0.1
$ws.Cells.Item(9, 10).Value = 0.1588235294117647
$ws.Cells.Item(9, 15).Value = 0.01764705882352941
$ws.Cells.Item(9, 17).Value = 0.1647058823529412
$ws.Cells.Item(9, 18).Value = 0.09411764705882353
$ws.Cells.Item(9, 19).Value = 0.4176470588235294
$ws.Cells.Item(10, 2).Value = 0.0687789799072643
$ws.Cells.Item(10, 4).Value = 0.02163833075734158
$ws.Cells.Item(10, 5).Value = 0.0007727975270479134
$ws.Cells.Item(10, 6).Value = 0.08114374034003091
$ws.Cells.Item(10, 10).Value = 0.1352395672333848
$ws.Cells.Item(10, 15).Value = 0.02782071097372488
$ws.Cells.Item(10, 17).Value = 0.2187017001545595
$ws.Cells.Item(10, 18).Value = 0.1004636785162288
$ws.Cells.Item(10, 19).Value = 0.3454404945904173
$ws.Cells.Item(11, 7).Value = 0.1628664495114006
$ws.Cells.Item(11, 10).Value = 0.07817589576547231
$ws.Cells.Item(11, 11).Value = 0.1954397394136808
$ws.Cells.Item(11, 12).Value = 0.5374592833876222
$ws.Cells.Item(11, 19).Value = 0.02605863192182411
$ws.Cells.Item(12, 7).Value = 0.7267441860465116
$ws.Cells.Item(12, 10).Value = 0.2151162790697674
$ws.Cells.Item(12, 11).Value = 0.01162790697674419
$ws.Cells.Item(12, 12).Value = 0.01162790697674419
$ws.Cells.Item(12, 19).Value = 0.03488372093023256
$ws.Cells.Item(13, 6).Value = 0.02083333333333333
$ws.Cells.Item(13, 7).Value = 0.7291666666666666
$ws.Cells.Item(13, 19).Value = 0.08333333333333333
$ws.Cells.Item(15, 6).Value = 0.02164502164502164
$ws.Cells.Item(15, 8).Value = 0.1731601731601732
$ws.Cells.Item(15, 9).Value = 0.06926406926406926
$ws.Cells.Item(15, 10).Value = 0.341991341991342
$ws.Cells.Item(15, 11).Value = 0.05194805194805195
$ws.Cells.Item(15, 15).Value = 0.03463203463203463
$ws.Cells.Item(15, 19).Value = 0.3073593073593073
$ws.Cells.Item(16, 6).Value = 0.007092198581560284
$ws.Cells.Item(16, 8).Value = 0.2269503546099291
$ws.Cells.Item(16, 9).Value = 0.07801418439716312
$ws.Cells.Item(16, 10).Value = 0.3687943262411347
$ws.Cells.Item(16, 11).Value = 0.1063829787234043
$ws.Cells.Item(16, 13).Value = 0.04964539007092199
$ws.Cells.Item(16, 15).Value = 0.05673758865248227
$ws.Cells.Item(16, 19).Value = 0.1063829787234043
$ws.Cells.Item(17, 6).Value = 0.006172839506172839
$ws.Cells.Item(17, 8).Value = 0.2098765432098765
$ws.Cells.Item(17, 9).Value = 0.06995884773662552
$ws.Cells.Item(17, 10).Value = 0.4259259259259259
$ws.Cells.Item(17, 11).Value = 0.08024691358024691
$ws.Cells.Item(17, 13).Value = 0.0308641975308642
$ws.Cells.Item(17, 15).Value = 0.07407407407407407
$ws.Cells.Item(17, 19).Value = 0.102880658436214
$ws.Cells.Item(18, 6).Value = 0.0128755364806867
$ws.Cells.Item(18, 8).Value = 0.2060085836909871
$ws.Cells.Item(18, 9).Value = 0.07296137339055794
$ws.Cells.Item(18, 10).Value = 0.4377682403433477
$ws.Cells.Item(18, 11).Value = 0.1072961373390558
$ws.Cells.Item(18, 13).Value = 0.02145922746781116
$ws.Cells.Item(18, 15).Value = 0.05150214592274678
$ws.Cells.Item(18, 19).Value = 0.09012875536480687
$ws.Cells.Item(19, 6).Value = 0.01798855273916599
$ws.Cells.Item(19, 8).Value = 0.2150449713818479
$ws.Cells.Item(19, 9).Value = 0.07277187244480784
$ws.Cells.Item(19, 10).Value = 0.3810302534750613
$ws.Cells.Item(19, 11).Value = 0.1226492232215863
$ws.Cells.Item(19, 13).Value = 0.01635322976287817
$ws.Cells.Item(19, 15).Value = 0.0678659035159444
$ws.Cells.Item(19, 19).Value = 0.1062959934587081
